$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 92, pushing existing rows 92-198 down to 93-199.
$ws.Rows.Item(92).Insert()

# Populate the newly-inserted row 92 with the new weekly price record.
$ws.Cells.Item(92, 1).Value = 8
$ws.Cells.Item(92, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(92, 3).Value = "Coquimbo"
$ws.Cells.Item(92, 4).Value = 45117
$ws.Cells.Item(92, 5).Value = 4
$ws.Cells.Item(92, 6).Value = 100112052
$ws.Cells.Item(92, 7).Value = "Albahaca"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 1200
$ws.Cells.Item(92, 11).Value = 3000
$ws.Cells.Item(92, 12).Value = 3500
$ws.Cells.Item(92, 13).Value = 3250
$ws.Cells.Item(92, 14).Value = "$/paquete"
$ws.Cells.Item(92, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(92, 16).Value = 3250
$ws.Cells.Item(92, 17).Value = 1
$ws.Cells.Item(92, 18).Value = "Hortaliza"
